# Generate Report for Handoff
#
# Re-running the handoff-report generation stamps a fresh "Latest Handoff
# Datetime" on every row that is currently sitting in a "Ready for
# handoff" (or failed-transform) state, for each localized-language sheet.
# Rows 7, 10, 11, 12, 13, 14, 15, 16 (column D) all collapse onto the new
# timestamp produced by this run, for both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$rows = @(7, 10, 11, 12, 13, 14, 15, 16)

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $ws_zhcn.Range("D$r").Value = "2016-03-11 02:36:09"
}

$ws_dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $ws_dede.Range("D$r").Value = "2016-03-11 02:36:17"
}
